$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# New row 60 - JP Flood general policy 'A2' with deductible
$ws.Range("B60").Value = "fm56"

# New row 61 - JP Flood all general policies unit test (T1 only)
$ws.Range("B61").Value = "fm57"

$ws.Range("C60").Value = "JP Flood general policy 'A2' with deductible"
$ws.Range("C61").Value = "JP Flood all general policies unit test (T1 only)"

$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 28
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = "complete"
$ws.Range("I60").Value = "complete"

$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 28
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = "complete"
$ws.Range("I61").Value = "complete"

# Update the existing row 59 description text (new naming for the "all step policies" test)
$ws.Range("C59").Value = "JP Flood all step policies unit test (T1,2,3,5)"
$ws.Range("E59").Value = "27,28,29,14"

# Copy the style from row 59 (B/C/H/I columns) down to the two new rows
$ws.Range("B59").Copy()
$ws.Range("B60:B61").PasteSpecial(-4122)
$ws.Range("C59").Copy()
$ws.Range("C60:C61").PasteSpecial(-4122)
$ws.Range("H59:I59").Copy()
$ws.Range("H60:I61").PasteSpecial(-4122)

$ws.Range("D61").Select()
